# Refresh the CV's "last rendered" metadata: the source R Markdown was
# re-knit on 2016-04-17 (was 2016-02-23), which bumps the derived age,
# the rmarkdown package version, and the recorded R sessionInfo() values
# that got captured in the "Documenting File Creation" section.

$d = $word.ActiveDocument

function Replace-ExactText($doc, [string]$old, [string]$new) {
    # Locate the run's exact text with Find (no in-place replacement text
    # supplied here) and then assign Range.Text directly - doing the
    # substitution through Find.Execute's own Replace argument runs the
    # new text through Word's AutoFormat "smart quotes" pass, which would
    # mangle the literal straight quotes (") that appear in several of
    # these strings.
    $rng = $doc.Content
    $found = $rng.Find.Execute(
        $old, $true, $false, $false, $false, $false,
        $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Replace-ExactText: text not found -> $old"
    }
    $rng.Text = $new
}

# Age at top of "1. About Me" table: 31y 4m 1d -> 31y 5m 27d
Replace-ExactText $d "31y 4m 1d" "31y 5m 27d"

# rmarkdown package version bump
Replace-ExactText $d "package version: 0.9.2" "package version: 0.9.5.9"

# "File latest updated date" bullet
Replace-ExactText $d "File latest updated date: 2016-02-23" "File latest updated date: 2016-04-17"

# Sys.time() "setting value" line (straight double quotes must stay straight)
Replace-ExactText $d "[1] ""2016-02-23 03:17:39 EST"" setting value" "[1] ""2016-04-17 22:20:39 EDT"" setting value"

# sessionInfo() date field
Replace-ExactText $d "date 2016-02-23" "date 2016-04-17"

# sessionInfo() sysname/nodename line - nodename changed from
# "rstudio-scibrokes" to "rstudio-scibrokes3"
$sysnameOld = "sysname release ""Linux"" ""3.10.0-229.20.1.el7.x86_64"" version nodename " + `
    """#1 SMP Tue Nov 3 19:10:07 UTC 2015"" ""rstudio-scibrokes"" machine login ""x86_64"" " + `
    """unknown"" user effective_user ""ryoeng"" ""ryoeng"""
$sysnameNew = "sysname release ""Linux"" ""3.10.0-229.20.1.el7.x86_64"" version nodename " + `
    """#1 SMP Tue Nov 3 19:10:07 UTC 2015"" ""rstudio-scibrokes3"" machine login ""x86_64"" " + `
    """unknown"" user effective_user ""ryoeng"" ""ryoeng"""
Replace-ExactText $d $sysnameOld $sysnameNew
